$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2025-02-02"
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 2
$ws1.Range("F2").Value = 2
$ws1.Range("G2").Value = 4
$ws1.Range("H2").Value = 8

$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "2025-02-09"
$ws1.Range("E3").Value = 3
$ws1.Range("F3").Value = 4
$ws1.Range("G3").Value = 6
$ws1.Range("H3").Value = 9

$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2025-02-16"
$ws1.Range("G4").Value = 8
$ws1.Range("H4").Value = 11

$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2025-02-23"

$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2025-03-02"
$ws1.Range("F6").Value = 7
$ws1.Range("H6").Value = 14

$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2025-03-09"
$ws1.Range("E7").Value = 6
$ws1.Range("G7").Value = 11
$ws1.Range("H7").Value = 16

$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2025-03-16"
$ws1.Range("G8").Value = 10
$ws1.Range("H8").Value = 17

$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "2025-03-23"
$ws1.Range("E9").Value = 6
$ws1.Range("H9").Value = 18

$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2025-03-30"
$ws1.Range("G10").Value = 10
$ws1.Range("H10").Value = 17

$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2025-04-06"
$ws1.Range("E11").Value = 6
$ws1.Range("G11").Value = 11
$ws1.Range("H11").Value = 19

$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2025-04-13"
$ws1.Range("F12").Value = 7
$ws1.Range("G12").Value = 11
$ws1.Range("H12").Value = 18

$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2025-04-20"
$ws1.Range("E13").Value = 6
$ws1.Range("F13").Value = 7
$ws1.Range("G13").Value = 11
$ws1.Range("H13").Value = 17

$ws1.Range("B14").NumberFormat = "@"
$ws1.Range("B14").Value = "2025-04-27"
$ws1.Range("G14").Value = 11

$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "2025-05-04"
$ws1.Range("E15").Value = 6
$ws1.Range("F15").Value = 7
$ws1.Range("H15").Value = 17

$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "2025-05-11"
$ws1.Range("E16").Value = 6
$ws1.Range("F16").Value = 7
$ws1.Range("G16").Value = 11
$ws1.Range("H16").Value = 18

$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "2025-05-18"
$ws1.Range("E17").Value = 6
$ws1.Range("F17").Value = 6
$ws1.Range("G17").Value = 11
$ws1.Range("H17").Value = 18

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "2022-12-25 to 2025-01-26"
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "14"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "7"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "3"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "0"
$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-02-02"
